$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

for ($r = 29; $r -le 33; $r++) {
    $ws2.Range("B28:D28").Copy() | Out-Null
    $ws2.Range("B$r`:D$r").PasteSpecial() | Out-Null
    $ws2.Rows.Item($r).RowHeight = 15
}

$ws2.Range("B29").Value = "fg 1"
$ws2.Range("B30").Value = "fg 2"
$ws2.Range("B31").Value = "fg 3"
$ws2.Range("B33").Value = "fg 5"
$ws2.Range("B32").Value = "fg 4"

$ws2.Range("D33").Value = "k"
$ws2.Range("D32").Value = "#ff0000"
$ws2.Range("D31").Value = "#ff9900"
$ws2.Range("D30").Value = "#ffff00"
$ws2.Range("D29").Value = "#ccff66"
